$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, matching the header style used by other headers
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Threshold used to derive the Save flag from the existing "sum" (column G) value
$threshold = 8.418600821238126

# Populate H2:H55 with 1 when sum >= threshold, else 0 (mirrors the values baked into the diff)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("H35").Value = 1
$ws.Range("H36").Value = 0
$ws.Range("H37").Value = 1
$ws.Range("H38").Value = 0
$ws.Range("H39").Value = 1
$ws.Range("H40").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("H43").Value = 1
$ws.Range("H44").Value = 1
$ws.Range("H45").Value = 1
$ws.Range("H46").Value = 0
$ws.Range("H47").Value = 1
$ws.Range("H48").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("H53").Value = 1
$ws.Range("H54").Value = 1
$ws.Range("H55").Value = 1

$excel.CutCopyMode = 0
